# Updated symbol list on Fri Dec 16 11:57:13 UTC 2022 with GitHub Actions
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Price (column D) refreshes - values are stored as text, so a leading
# apostrophe is used to force text entry (matches how these cells were
# already stored, avoiding Excel's automatic number inference).
$priceUpdates = @{
    "D2"  = "251.51"
    "D3"  = "23.75"
    "D4"  = "5.995"
    "D5"  = "0.05964"
    "D6"  = "3.423"
    "D8"  = "1.321"
    "D9"  = "0.7983"
    "D10" = "0.1489"
    "D11" = "0.07866"
    "D12" = "0.03350"
    "D14" = "0.09277"
    "D15" = "3.558"
    "D16" = "0.001675"
    "D18" = "0.0006082"
    "D19" = "0.006221"
    "D21" = "0.001068"
    "D22" = "0.0001501"
    "D23" = "3.684"
    "D25" = "0.3308"
    "D27" = "0.0006477"
    "D40" = "0.04445"
    "D41" = "0.007005"
    "D44" = "0.009240"
    "D46" = "0.00005882"
    "D49" = "0.09782"
}

foreach ($addr in $priceUpdates.Keys) {
    $ws.Range($addr).Value = "'" + $priceUpdates[$addr]
}

# Rows 42 and 43 swapped rank order: BKEXToken now sits above CEJI.
$ws.Range("B42").Value = "BKEXToken"
$ws.Range("C42").Value = "https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk"
$ws.Range("D42").Value = "'0.1069"
$ws.Range("E42").Value = "41BKEXTokenBKK"

$ws.Range("B43").Value = "CEJI"
$ws.Range("C43").Value = "https://coinranking.com/coin/SbKjCVJCh+ceji-ceji"
$ws.Range("D43").Value = "'0.003151"
$ws.Range("E43").Value = "42CEJICEJI"
